$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: date, hours spent, and description of work done
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 44312
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = "Сделал простой сервис, раскидал по папкам файлы."

# Move the active selection as recorded after the edit
$ws.Range("C4").Select()
